$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 11629149
$ws.Range("I43").Value = 62501852
$ws.Range("J43").Value = 1102.0571
$ws.Range("K43").Value = 62501852
$ws.Range("L43").Value = 1102.0571
$ws.Range("M43").Value = -62501783
$ws.Range("N43").Value = -1240.0571
$ws.Range("H132").Value = 9808366
$ws.Range("I132").Value = 19235184
$ws.Range("J132").Value = 4475.96
$ws.Range("K132").Value = 57705552
$ws.Range("L132").Value = 13427.88
$ws.Range("M132").Value = -57703022
$ws.Range("N132").Value = -18487.88

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 63787.375
$ws.Range("I45").Value = 125812.375
$ws.Range("J45").Value = 1762.375
$ws.Range("K45").Value = 125812.375
$ws.Range("L45").Value = 1762.375
$ws.Range("M45").Value = -125435.375
$ws.Range("N45").Value = -2516.375
$ws.Range("H110").Value = 2209.6667
$ws.Range("I110").Value = 871.6
$ws.Range("J110").Value = 8900
$ws.Range("K110").Value = 871.6
$ws.Range("L110").Value = 8900
$ws.Range("M110").Value = 1173.4
$ws.Range("N110").Value = -12990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3245.3333
$ws.Range("I20").Value = 3164.25
$ws.Range("J20").Value = 3569.6667
$ws.Range("K20").Value = 3164.25
$ws.Range("L20").Value = 3569.6667
$ws.Range("M20").Value = -2917.25
$ws.Range("N20").Value = -4063.6667
$ws.Range("H69").Value = 15888.5
$ws.Range("J69").Value = 15888.5
$ws.Range("L69").Value = 15888.5
$ws.Range("N69").Value = -17510.5
$ws.Range("H72").Value = 15888.5
$ws.Range("J72").Value = 15888.5
$ws.Range("L72").Value = 47665.5
$ws.Range("N72").Value = -55777.5
$ws.Range("H80").Value = 4052.303
$ws.Range("I80").Value = 855.55554
$ws.Range("J80").Value = 5251.0835
$ws.Range("K80").Value = 855.55554
$ws.Range("L80").Value = 5251.0835
$ws.Range("M80").Value = 142.44446
$ws.Range("N80").Value = -7247.0835
$ws.Range("H83").Value = 4052.303
$ws.Range("I83").Value = 855.55554
$ws.Range("J83").Value = 5251.0835
$ws.Range("K83").Value = 4277.7777
$ws.Range("L83").Value = 26255.4175
$ws.Range("M83").Value = 714.2223000000004
$ws.Range("N83").Value = -36239.4175
$ws.Range("H107").Value = 45455548
$ws.Range("I107").Value = 45455548
$ws.Range("K107").Value = 45455548
$ws.Range("M107").Value = -45453628
$ws.Range("H134").Value = 4832110
$ws.Range("I134").Value = 1081.15
$ws.Range("J134").Value = 37038972
$ws.Range("K134").Value = 3243.45
$ws.Range("L134").Value = 111116916
$ws.Range("M134").Value = -708.4500000000003
$ws.Range("N134").Value = -111121986

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 25055652
$ws.Range("I86").Value = 45504910
$ws.Range("J86").Value = 62112
$ws.Range("K86").Value = 45504910
$ws.Range("L86").Value = 62112
$ws.Range("M86").Value = -45503787
$ws.Range("N86").Value = -64358
$ws.Range("H89").Value = 25055652
$ws.Range("I89").Value = 45504910
$ws.Range("J89").Value = 62112
$ws.Range("K89").Value = 227524550
$ws.Range("L89").Value = 310560
$ws.Range("M89").Value = -227518934
$ws.Range("N89").Value = -321792
$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164
$ws.Range("H132").Value = 37039440
$ws.Range("I132").Value = 2019.8
$ws.Range("J132").Value = 83336210
$ws.Range("K132").Value = 6059.4
$ws.Range("L132").Value = 250008630
$ws.Range("M132").Value = -3529.4
$ws.Range("N132").Value = -250013690
$ws.Range("H134").Value = 985.62964
$ws.Range("I134").Value = 1400.9231
$ws.Range("J134").Value = 600
$ws.Range("K134").Value = 4202.7693
$ws.Range("L134").Value = 1800
$ws.Range("M134").Value = -1667.7693
$ws.Range("N134").Value = -6870
$ws.Range("H140").Value = 52350
$ws.Range("J140").Value = 52350
$ws.Range("L140").Value = 52350
$ws.Range("N140").Value = -62710

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 83.71429000000001
$ws.Range("I6").Value = 67.07692
$ws.Range("K6").Value = 201.23076
$ws.Range("M6").Value = -88.23076
$ws.Range("H37").Value = 222780.6
$ws.Range("J37").Value = 222780.6
$ws.Range("L37").Value = 668341.8
$ws.Range("N37").Value = -668565.8
$ws.Range("H122").Value = 7250389
$ws.Range("I122").Value = 29412040
$ws.Range("J122").Value = 5233.731
$ws.Range("K122").Value = 264708360
$ws.Range("L122").Value = 47103.579
$ws.Range("M122").Value = -264705910
$ws.Range("N122").Value = -52003.579
$ws.Range("H131").Value = 848.73
$ws.Range("J131").Value = 861.5876500000001
$ws.Range("L131").Value = 2584.76295
$ws.Range("N131").Value = -12664.76295

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2558.2856
$ws.Range("I132").Value = 1983.8823
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 5951.6469
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -3421.6469
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H100").Value = 2389.65
$ws.Range("I100").Value = 2260.3
$ws.Range("J100").Value = 2432.7666
$ws.Range("K100").Value = 2260.3
$ws.Range("L100").Value = 2432.7666
$ws.Range("M100").Value = -1719.3
$ws.Range("N100").Value = -3514.7666
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 45716624
$ws.Range("I136").Value = 6496017.5
$ws.Range("J136").Value = 333334400
$ws.Range("K136").Value = 19488052.5
$ws.Range("L136").Value = 1000003200
$ws.Range("M136").Value = -19485502.5
$ws.Range("N136").Value = -1000008300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 80000750
$ws.Range("J86").Value = 80000750
$ws.Range("L86").Value = 80000750
$ws.Range("N86").Value = -80002996
$ws.Range("H89").Value = 80000750
$ws.Range("J89").Value = 80000750
$ws.Range("L89").Value = 400003750
$ws.Range("N89").Value = -400014982
$ws.Range("H113").Value = 652.75
$ws.Range("I113").Value = 711.7143
$ws.Range("J113").Value = 240
$ws.Range("K113").Value = 2135.1429
$ws.Range("L113").Value = 720
$ws.Range("M113").Value = 34.85710000000017
$ws.Range("N113").Value = -5060
$ws.Range("H132").Value = 45675.668
$ws.Range("I132").Value = 53910.9
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 161732.7
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -159202.7
$ws.Range("N132").Value = -18558.5
